# Manage.xlsx - "make a report excel" edit
# Adds new Employee rows (Seoyoung, Ayoung) and new Project rows
# (첫 프로젝트 / 알바 충원) to the employee/project report workbook.
#
# Note: values that look like numbers or dates (e.g. "030416", "0",
# "2024-08-03") are entered with a leading apostrophe so Excel stores
# them as literal text (shared strings) instead of silently converting
# them to numbers / date serials, which would corrupt data such as the
# leading zero in "030416".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Employee": add row 3 (Seoyoung) and row 4 (Ayoung)
# ---------------------------------------------------------------
$wsEmployee = $wb.Worksheets.Item("Employee")

# Row 3 - Seoyoung
$wsEmployee.Range("A3").Value = "Seoyoung"
$wsEmployee.Range("B3").Value = "'030416"
$wsEmployee.Range("C3").Value = "010-5718-4778"
$wsEmployee.Range("D3").Value = "seoyoung.you@gmail.com"
$wsEmployee.Range("E3").Value = "'0"
$wsEmployee.Range("F3").Value = "인턴"
$wsEmployee.Range("G3").Value = "이메일이 맞는 지 모르겠음"

# Row 4 - Ayoung
$wsEmployee.Range("A4").Value = "Ayoung"
$wsEmployee.Range("B4").Value = "'030416"
$wsEmployee.Range("C4").Value = "010-5719-4778"
$wsEmployee.Range("D4").Value = "'"
$wsEmployee.Range("E4").Value = "'0"
$wsEmployee.Range("F4").Value = "인턴"
$wsEmployee.Range("G4").Value = "이메일 부재"

# ---------------------------------------------------------------
# Sheet "Project": add row 2 (첫 프로젝트) and row 3 (알바 충원)
# ---------------------------------------------------------------
$wsProject = $wb.Worksheets.Item("Project")

# Row 2 - 첫 프로젝트
$wsProject.Range("A2").Value = "첫 프로젝트"
$wsProject.Range("B2").Value = "'2024-08-03"
$wsProject.Range("C2").Value = "'2024-08-22"
$wsProject.Range("D2").Value = "카나타"
$wsProject.Range("E2").Value = "Taeyoung"
$wsProject.Range("F2").Value = "Seoyoung"

# Row 3 - 알바 충원
$wsProject.Range("A3").Value = "알바 충원"
$wsProject.Range("B3").Value = "'2024-09-03"
$wsProject.Range("C3").Value = "'2024-09-04"
$wsProject.Range("D3").Value = "카나타"
$wsProject.Range("E3").Value = "Taeyoung"
$wsProject.Range("F3").Value = "Seoyoung"
$wsProject.Range("G3").Value = "Ayoung"
